# Update the YouTube URL stored in B2 and leave the selection on that cell,
# matching the author's re-upload of vurl_list.xlsx (commit: "Add files via
# upload").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "https://www.youtube.com/watch?v=PfTKqh_aqJE"

$ws.Activate()
$ws.Range("B2").Select()
